$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "战斗机"
$ws.Range("B5").Value = "轰炸机"
$ws.Range("B3").Select()
